$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13-17 replicate the same order (row 12's data) five more times.
# The source text for N/O/P (PGP key, encrypted messages, plaintext notes) is
# expressed as single-quoted here-strings (no interpolation) with one extra
# leading apostrophe, because Excel silently consumes the first apostrophe of a
# cell value as a "treat as text" quote-prefix; doubling it keeps one literal
# apostrophe in the stored string, matching the source data.
$pgpKey = @'
''-----BEGIN PGP PUBLIC KEY BLOCK-----
mQINBFnEIiEBEADMrm4vb5QBopQe5Jq7fM/B0qoiICt0+H1NzbiY7vyG8a4LmOMB
IIdfSjlfhdRC3snzL1t+cYZR13/msqlIRCjHlamGgHCoOaJftd0f+lUdopVWKF+c
RGMaEfnN4ZsS2VFTZnlGvIP6roZtNazD7XoGprGWqht30WUw0gkzgn3540zYJmMO
ijz/i9BFYhvAtMo9ABzdxMkK6kdNr4d7QiuhETw4n8NhgPro7BvciQxVDwvDKUAE
49NUPpXTgnEugemQS4E3hXZ6O5YbBVsSZJCKOd20tkZBAYUKYd834TEcJ7V8I7L/
xOeAiwQYHCZiIYpdO6BKpPl2XKHUyCcHp47jgN+jJpsYSyJrlN3llt1mc1bl42ei
smql9y4PiLXo4uhQoaJrDNc/MeD7XMFYrQXpuckVNXeu/EhwngQZrDuznPv4ncOz
qfc6Lw73vAioUli3Rpfv3LVFmtJbXHJHcSo6/5+bLLznT113XrgBoZaAaVtAVaiV
1kDX5PE0b0mzOw7mRNzw9Q+gVmSiPbO//JOpJfqIaZuJ51umjOyNlgQ1VQ+1TVY7
s75IPnnZb6/J3AzZjh3qB4JrK3E0WydgeM9vrqKzOv5G+Ro1Wwb4QF4O1GrF9e7u
Fa+hzyMO9ASRINISR3+5RCvmghxrYyBRfDoQPwlYRryA4tkZ2BtU+AeEMQARAQAB
tCBCYXJiYXJhIDxiNHJiNHI0QHByb3Rvbm1haWwuY29tPokCOAQTAQIAIgUCWcQi
IQIbAwYLCQgHAwIGFQgCCQoLBBYCAwECHgECF4AACgkQSw6XQ6wu7xhIehAAtxfI
/WFCwZzyyailymwphJm8EYW0RU157r5mm/NB5E4W4SUlibnLFk+ThImAfhn2kmW8
2a0O39jzgWfqWhwCxsbSrDsg3eC9X5jdkZIid9e1C5oZVGVDeA3ScowxApFAdEZp
nskR05uw72zK8TgCLlunNwAx+BBk9pzl3eys2Xgwj/+fQhOzFhIVvKUKSpA9MZ5c
WaE8w+BD0SPHNof+btwoSzy1FeHhAZ2umaesXIsAJI+jGxPQi7h/82k9ZJ5jitfw
gfSKIYEptIS2RO4L0hDVYb0Spy1q5pfYDYobh7QaEVxQ1emijixZuPwhRJZUOaQj
+6sC1yWDifKE4Td135tn5MBQkjFXziTWQRVP8aJpZYeAq4XO3g4miBc/X3f+HTrz
3c1NCJ05SjJ/7kufzup2rmxHhXF9qIRKtr84u/PJrV0KftSdxbuV6S+l7M2c9BnP
mjrJse4g/7pRowzQSljzXAlSy8B9sIONM07HhRJloqbw022u+nlpFLoR9gzkgrDp
1zVm7wzjk7X78p/mWLqf9iqqovtHbx61Uja4Rc6I3GxKBNGruifVv0R5CL1pQe6J
VnUZ3OpWKa4xHX0/EOywacvCFgfl4UlZWX8fBbZZWvsVKhzXpCSD2WqKjAYQrRi+
KwG8AgA+WOKLJOzwUshK4a0OvLiGVugPK6Vj0165Ag0EWcQiIQEQAMoPxzzljR7Y
xW1/aQ5udRQ1CkibW8FAOR4GuyTJmR2lLXU5Z4k+a6E4gIyM2bZs27LKG1Xc8dny
e7VuaIMkntgTnHIn8qqiU2s2nTEXDGg4D9zDLpsQWv+ABs5oByEK46IAnI62IUSY
Lb2pFTMvEZ+FR2Zk3foILg4znoyXAlzsHIhmyuvm6j2qbOV0zW7QnUg1ZmngPrhM
URgFU/+qN6/DJKcFXMAAEjxTaby2HWGxWPt5+vfTJAMbQVZIjdvSlQ+tEf0XWT7g
3qxx+DU5/x4uVh4wW9tOiiK/okE0rEoaIWwltu0A8c2u80kfiSv9zGOQ4xNLfPQr
YrlNKNPKIa2qgiHHIf7rt+Myxn9E6tFUJ5ZtB95i28DydYQTZjcDaax49zo9j2q+
Emj7NTZ8ea9CEGL358jLUOcCpwiz57pSlTUvowu8ckFfTNm7xLjSddbN24KMdPvY
lJHMc8CwKi9R0PSJJzrPdj5enFVPazGgMfH8xnR6FWCvPVhQDJRZQViC8+y1NCVo
+Vc/FWwDp5JsagH63lNj+d8J6VvgZ6dKTfR3TBj27586N2WuD0oox8etO7BqRmX0
f/39nOUm2wFnhSRi6DcucDeO+QReq+7T7W+zXQc/mUb8xFgUKIVAEfRQCsjXbLNk
/eS5MdEWAuMMuomAvNkPSu+IYRhMjVZdABEBAAGJAh8EGAECAAkFAlnEIiECGwwA
CgkQSw6XQ6wu7xgoZhAAu9qsEM1czZDD6qqwGRv7zMAnTireGvRmE60Y1NwSVowq
Wv2YNDzQLr1FcCcdZAuLm85OOl5TMaeAncjpdvWiSKlFtEyn4qAvWWxLPOB+Iue6
3y65h7+/ZzWjRUEPamTJjZf/87DY5TZiSWSzMMvAyugGGpPwnsu0KfWjinUX2Nfm
EVcj9NWkpNkna+sS7pq7tDN1q+8Lm3dIOv9CEnxxI9jgi1RUiGELthq46lHadWux
x0UV1yvn8QtAajRAuqZYgwi7oLYsa0n0nEgjEq5H2Mn5j1m5LwlFbrRcZb2v+5FB
fgcxs2VW3nbEOUDtlxCD3ySE2jHksxjRoJTuQRPPXDLAphoqV3OrJnaSl6UzWgaB
BjX++NXdb5TVAzDNRDQicSbP7SWdtTQD/G5dR8Py19DQF97Yiolw+WYbwSBde3H5
e2ecnJGObVWehfZ0H4hb7ekgWH3Gw8RFExd6ydSt0Beaf/B8ZPvw1LZ8ZN2qaeFt
EOgFxUn59D1uzxxXJNw9zop1mUfELfcrKD4vgxIZXrLWw1n9i8z6Au1eud+znOZM
rBb0Fx0i0pTVF/QZ98UQ0w4+cyGLS9VtYZ0QZFY3e4SbLfEMNoidPc/JTNrVXzxF
gW5Qju7JUghRY7ti5QcNXZDu1WFzwYaTUrD4ZBVKf4azSWOt99i0EzBmDSvKdiM=
=qzG7
-----END PGP PUBLIC KEY BLOCK-----
'@

$pgpMessages = @'
''-----BEGIN PGP MESSAGE-----
hQIMA1XtIn7QB2h8ARAAoRfX/zZ5DShXY7dBkk428DjLovHsBgTvzuKkQerrBI1y
NJQSXMdX+71vJFaesjojHamv7kcBlBDYkZc7CETOsekDCJPC7a1OBGFZ20AV39Ep
uni7vbLrR0TJ/d0o+t66bTf7NJ2ZeRlertIgiJjUtVwobgOZls/uwpVOEWRZ3g76
+9M5aJtzCxnQN2gKdVqVvgnNNGCDE1aZw9TO1AXyK/Z51vWYrIPbQ2fcWHhE44Ya
6YslLxTHEqqN0OIBR5U/VxZ52OegmOcOL9YzUp9Coh+fNyczRvCrM1pjGbODffkX
j4JaNPaUxbtbnZfGg0B6U27U1GN3PC793F5c7RoicM5RssIUC4tBP+/YatTnIRkv
x//2BZ7MXdIn0OVUkdDKZEadJMzXIpnMTlVfFKDbBY/j+ghf7wBFTjzVk6teJEj2
Azh54zkLLwldWHxAX/uSotNQ4/zUqdi3PRIqECp4wifx2bNLqXE3Depu2LVtdfcV
CLLxuyPv0T7MNqdsRURqqvHDglbU0WTJMFPDhz5DAFNpNNor6PjUxPaP0bOvRxiN
jIh4qGjiB6sm+38bYBa8SVH4BaqXUnmwoWTGzrAuatxfYqIAdx0BRzFqM5wbS4Ed
KEE9Uonr61I/ZntghQ0f1zqO5kCEraY993z6O3ScEZIXPVaOXfY3TojHW7ApNlqF
AgwDVjkqg7iqkkQBEACDd5+NgO6bnOgzEq0nKHXxnNv5zwG6f8Ts0BnXY34RZANn
3UmRYziOV8naxlKZbHulix33NC+m/3ofaDDL/9UQzPSvUUwWWmkb5/3U70c66KCN
jXWLUyXgGAi1EObCgzA+8Nd2QwpVgJ5QXY/aTHTwJjoiL82z0gIuqZ/loVzlbyB+
S6L0QU7Tec1pDnAGu4VMKmaQ7l17xAN7AkxVRK/cxbMf9KIyf5m3W8d97LBuCtdW
yfqpKYDD/jBXQ/Ta1xQzj5kfKLpKVqTnaFC3gsMUzI9e9NCdsswwKJ38qAjAyTV+
xPTxJberF+P6q2ufsCshgrCoNig1Uu0ZMI4bcr6TLouzg05UIzf+Z7x23YQihknx
n3BVRyW7DWpgmuxavLdoyQgIqlnhWTm7LOel0XBMP8q8LKyoOaxYyb3QGlEuDcNE
HDtkqNtPDmBJwRg2UkCy/VsagqCUb+pl/mi81W/+g6fhMkdZGladTuI2U6e1fABC
BG0l2v68n+1GXBit7nwFP5uZhjjw4tptk7MQflSQFGsGWkh4yjEJqtd/tsF+hSJE
+A/AD00EJy5a/UGl0ZAdq71e1TQZr3NJQGFILTUvraVX8eVJzEBKmW0bHexJ9Doi
tMIihd896jx/lwbb5IFUNRUeVdST/VrxCz2DwnjF0HXqiAU98FWbTzkTD1vf0tLp
AQpXPDCcwmx0uhsBRFC8hcmgzKtjDA+fmxKVO96yzTbDgrzxyaotGL+QTOoNaEcC
h3/SmaEPdg5/0++gzrIzH2dvD8ScurTQ4o47ygL0swUHkIX4CY+Xwo4Hnh9MIJwi
AyEFIO1gLlq1yLMCw2qAmCfBU3lfZ+EONZdiMu5gHf/CTmUmV/zY2Ar7V7J7Orb3
AVzX3Sexx3AL4wscU9CDjywsEWQB9CEfHA88K9N5/ALTms1kecR0KP2XXvltiBLA
haQP0m8bv0SesVii/Kfb9v2HyChBa3NN18c6Ks3tN/aFnH4v5wuYIGDtvICTCG+0
A+UNZyAkbg/Ykr4X5wPI33aWGPL6viBpO+3i+5mir1PLBPnRcdmazIF54bSs2gur
lNce+RctppqOhpul529vlc42K1FEuzaTm4W9D1Z96m4PsMezjhMr1gMUEa598628
t9V699r540f5Qv93NBNe7UFgWJceaIWJzhvM55zKtXBLw7LU8+oxWQ64mpVRgoZv
BdjiRO6RoFCnTUGmCzuV9pKddZ7H7dD+VDovk6QI7mVG755+9Iw49B2o3doD6Ft2
NJFE3Jrgl4ZHZMY+boc3x97j4tK+bGBqJ5k2rJFtbBEMJO/ymB4EhgYD8VnXFftZ
j25V52NyBuCp75CsMhn3fhCtrioln9X4ivbk+FqQxpiVAxOyepjD5IJUQDW3xErw
SQAqQoUr+H9BMtWjEHfp00pmO/Uzq/TD+ohP0nsvL3Vh8gPYeYQXjbRfeg76a9v+
hLREyIPXsQSeFZ82XtVzkklv5UZqMmxJDrGI5o7Eduql1iO34qlgpfFyHSLrembi
B9xnITM9xHNJ7CrAXzLEjKLrynKPT0etW/GWBeKk4uSKX8MEZ4I=
=mO3k
-----END PGP MESSAGE-----
-----BEGIN PGP MESSAGE-----
hQIMA1Y5KoO4qpJEAQ//R26qF9UY/Sh0e0HyCsD+kvhUWxNWOdoD3YJSLZfb1bu1
u3J/v0nQUoplzMC/N1oxTf0ZBgQFUSRXT0GTzmQNZknwsZDPZMIPs9n9JTrAo5yq
GDF2CrUrQKJJz8pMZ0r2eaig7lUmrdN4/Ruve66rrF2vxwDa4j2kR+HTq3pfNoof
UwAS8ZhXneRYLCP+BNV41ea5eWHeqj8cVGJnW0KIm34qfCLJKrs4J3C8KjQQsAHA
SmYDLhu/gntZPN1v/+zxZ0fXltkzeN99HeXJmnyLRd36YmgST34Y+OfrACpx2PSn
aNBhCVFLvt3wN+AeSFOAZUIMAvwdNtvSOgeGq3FCgeHpE2R9I3626SssufVfDgpq
bwxHiG37wLF7GSlIo0wOcuqlmuKBa7b2q9gINefA2+FSpU6CrvNwNLCiO6yBDmAg
Ldv9bXKEw6ZbNhTGICQ3FMwkB9mDDdJ+qzCS2oVsC91nvvnB8huVEqNlA06URBEX
564x/cdBRiqCHaRNXii7eNNsrupTJJCEiY3cq8rewTH5Stc8Y/Ju2E8a7uHUJszR
GtxgSpvJlvXRm41eX/SJ1Jel2gIUdKmNX2Lq9ZU0Ff8HUPCZ5QJHR4rpf5yY6IJk
7SSqLMI7s8/UhUm9rSwLEhChJBZt8Gbg2UH9Swx/IIvEmNbcq44rIrjgSY9t0OPS
bwFoThLJ8SFWLEZ1owP2ws3kzWZOkRiSRtIBY/HBeTXziHpOsVp1Gd122DMGoUom
/FFe8hL2r9iJXdn9FxoJ3z/vsBnVbEZqNdz6ORG9rQDSEKVVbmT+0GM4IThQklYs
cTS7RKCq++tUWmJvn9chfQ==
=OUrC
-----END PGP MESSAGE-----
'@

$plainMessages = @'
''Order processed by vendor.
Order accepted by vendor.
'@

for ($r = 13; $r -le 17; $r++) {
    $ws.Range("A$r").Value = '33fae5b6'
    $ws.Range("B$r").Value = '2025-08-05 04:05:02'
    $ws.Range("C$r").Value = 'darkwangduck'
    $ws.Range("D$r").Value = 'Ben Eckart'
    $ws.Range("E$r").Value = '3873 Patterson Ave'
    $ws.Range("F$r").Value = 'Oakland'
    $ws.Range("G$r").Value = 'CA'
    $ws.Range("H$r").Value = '94619-2029'
    $ws.Range("I$r").Value = 'wax melts -  Punishers'
    $ws.Range("J$r").Value = '[Bulk option] 300g - Wax Skulls'
    $ws.Range("K$r").Value = 'UPS 2 Day'
    $ws.Range("M$r").Value = 'A0FAB2C4179732922AEEE4924B0E9743AC2EEF18'

    $ws.Range("N$r").Value = $pgpKey
    $ws.Range("N$r").Style = "Normal"
    $ws.Range("N$r").WrapText = $true
    $ws.Range("N$r").VerticalAlignment = -4160

    $ws.Range("O$r").Value = $pgpMessages
    $ws.Range("O$r").Style = "Normal"
    $ws.Range("O$r").WrapText = $true
    $ws.Range("O$r").VerticalAlignment = -4160

    $ws.Range("P$r").Value = $plainMessages
    $ws.Range("P$r").Style = "Normal"
    $ws.Range("P$r").WrapText = $true
    $ws.Range("P$r").VerticalAlignment = -4160

    $ws.Range("Q$r").Value = 2.334189202708
    $ws.Range("R$r").Value = 700.12
    $ws.Range("S$r").Value = 0.133382240155
    $ws.Range("T$r").Value = 40.01
    $ws.Range("U$r").Value = 2.467571442863
    $ws.Range("V$r").Value = 740.12

    # Restore default (non-custom) row height after the multi-line wrapped
    # text nudged Excel into recording an explicit row height.
    $ws.Rows.Item($r).AutoFit()
}
